$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.739.09'
$ws.Range('E2').Value = '  +4.60%  '
$ws.Range('D3').Value = '3.620.18'
$ws.Range('E3').Value = '  +3.75%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''629.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.31%  '
$ws.Range('D6').Value = '''158.54'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.78%  '
$ws.Range('D7').Value = '3.618.60'
$ws.Range('E7').Value = '  +3.74%  '
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('E9').Value = '  +3.83%  '
$ws.Range('D10').Value = '''0.149'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.09%  '
$ws.Range('D11').Value = '''7.37'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.83%  '
$ws.Range('E12').Value = '  +4.37%  '
$ws.Range('E13').Value = '  +5.43%  '
$ws.Range('D14').Value = '''33.37'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.14%  '
$ws.Range('D15').Value = '4.236.09'
$ws.Range('E15').Value = '  +3.49%  '
$ws.Range('D16').Value = '69.859.06'
$ws.Range('E16').Value = '  +4.60%  '
$ws.Range('D17').Value = '3.620.89'
$ws.Range('E17').Value = '  +2.84%  '
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').Value = '''6.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.18%  '
$ws.Range('E20').Value = '  +6.39%  '
$ws.Range('D21').Value = '''10.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +15.27%  '
$ws.Range('D22').Value = '''464.02'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.52%  '
$ws.Range('D23').Value = '''0.646'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.31%  '
$ws.Range('D24').Value = '''78.82'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.63%  '
$ws.Range('E25').Value = '  +13.98%  '
$ws.Range('D26').Value = '''10.72'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.08%  '
$ws.Range('D27').Value = '3.767.53'
$ws.Range('E27').Value = '  +3.34%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').Value = '''9.21'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +14.29%  '
$ws.Range('D30').Value = '''2.63'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.38%  '
$ws.Range('E31').Value = '  +10.86%  '
$ws.Range('E32').Value = '  +13.60%  '
$ws.Range('D33').Value = '''6.57'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.28%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').Value = '''1.96'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.95%  '
$ws.Range('D36').Value = '''26.59'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.09%  '
$ws.Range('D37').Value = '3.621.15'
$ws.Range('E37').Value = '  +3.73%  '
$ws.Range('D38').Value = '''8.47'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.26%  '
$ws.Range('D39').Value = '''2.43'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +14.69%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').Value = '''0.0926'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.94%  '
$ws.Range('D42').Value = '''179.46'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.00%  '
$ws.Range('D43').Value = '''1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').Value = '''5.69'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.38%  '
$ws.Range('D45').Value = '''32.50'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +21.51%  '
$ws.Range('D46').Value = '''0.915'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.89%  '
$ws.Range('D47').Value = '''1.37'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +14.91%  '
$ws.Range('D48').Value = '''46.22'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.47%  '
$ws.Range('D49').Value = '''2.76'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +13.22%  '
$ws.Range('D50').Value = '''7.82'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.70%  '
$ws.Range('B51').Value = 'Bittensor'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D51').Value = '''369.64'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +13.85%  '
